# Automatische test-sync: 2025-06-26 22:22:50
$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row 24 with the latest test mail ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A24").Value = "Wanneer zijn jullie open?"
$logs.Range("B24").Value = "mailmind.test@zohomail.eu"
$logs.Range("C24").Value = "Testmail #1: Wanneer zijn jullie open?"
$logs.Range("D24").Value = "Openingstijden / Locatie"
$logs.Range("E24").Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F24").Value = "2025-06-26 22:22:21"
$logs.Range("G24").Value = "Ja"
$logs.Range("H24").Value = "Nee"
$logs.Range("I24").Value = "Ja"

# Extend the conditional formatting ranges from row 23 to the new row 24
$logs.Range("D2:D23").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D24"))
$logs.Range("G2:G23").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G24"))
$logs.Range("H2:H23").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H24"))
$logs.Range("I2:I23").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I24"))

# --- Sheet "Dashboard": categories shift down by one row (wrapping),
#     reflecting the newly added "Openingstijden / Locatie" mail ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Openingstijden / Locatie"
$dash.Range("A4").Value = "Offerte / Prijsaanvraag"
$dash.Range("B4").Value = 2
$dash.Range("A5").Value = "Retour / Terugbetaling"
$dash.Range("A6").Value = "Productinformatie"
